# F-825-1032F CMD-J.L00.00.C aprevo TLIF-CA Fnl QC Insp.xlsx
# -------------------------------------------------------------------------
# Template update: the numbered "{{Item N}}" merge-field placeholders in the
# Cover Sheet's item grid (C3:D11) are re-filled column-major instead of
# row-major (C column = Item 1-9 top to bottom, D column = Item 10-18 top to
# bottom). In the process, the placeholder that used to read "{{Item 7}}"
# (with a space) is renamed to "{{Item7}}" (no space) so it can no longer be
# confused with any other numeric tool item text downstream.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cover Sheet")

# Column C: Item 1 .. Item 9 (row 9 uses the renamed "{{Item7}}" token)
$ws.Range("C3").Value  = "{{Item 1}}"
$ws.Range("C4").Value  = "{{Item 2}}"
$ws.Range("C5").Value  = "{{Item 3}}"
$ws.Range("C6").Value  = "{{Item 4}}"
$ws.Range("C7").Value  = "{{Item 5}}"
$ws.Range("C8").Value  = "{{Item 6}}"
$ws.Range("C9").Value  = "{{Item7}}"
$ws.Range("C10").Value = "{{Item 8}}"
$ws.Range("C11").Value = "{{Item 9}}"

# Column D: Item 10 .. Item 18
$ws.Range("D3").Value  = "{{Item 10}}"
$ws.Range("D4").Value  = "{{Item 11}}"
$ws.Range("D5").Value  = "{{Item 12}}"
$ws.Range("D6").Value  = "{{Item 13}}"
$ws.Range("D7").Value  = "{{Item 14}}"
$ws.Range("D8").Value  = "{{Item 15}}"
$ws.Range("D9").Value  = "{{Item 16}}"
$ws.Range("D10").Value = "{{Item 17}}"
$ws.Range("D11").Value = "{{Item 18}}"

# Leave the active selection on D8, matching the last cell touched in the
# authored session.
$ws.Range("D8").Select()
